# pabi_account_report - Advance Payment report template
# Rearrange / add columns in the header row (row 8):
#   - Move "Journal Number" / "Description" from L:M to G:H
#   - Insert new "Document Date" / "Posting Date" labels at I:J
#   - Move "PO Number" / "Contract Number" from G:H to N:O
#   - Everything from (old) N:Y shifts right by two columns, to P:AA
#
# The column-level width/style metadata is NOT shifted (it stays tied to the
# column index), only the cell text moves - so we only touch row 8 values,
# then separately widen the columns that now hold longer header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the old header row 8 text, column by column, before overwriting ---
# (use Value2 for reading - the bare .Value property getter doesn't resolve
# through this COM shim the way .Value2/.Text do)
$old = @{}
foreach ($col in @("G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y")) {
    $old[$col] = $ws.Range($col + "8").Value2
}

# --- new labels that don't exist yet in the sheet ---
$docDate = "Document Date"
$postingDate = "Posting Date"

# --- write the new header row 8 layout ---
$ws.Range("G8").Value = $old["L"]   # Journal Number
$ws.Range("H8").Value = $old["M"]   # Description
$ws.Range("I8").Value = $docDate    # Document Date (new)
$ws.Range("J8").Value = $postingDate # Posting Date (new)
$ws.Range("K8").Value = $old["I"]   # Supplier Invoice Number
$ws.Range("L8").Value = $old["J"]   # Supplier Invoice Date
$ws.Range("M8").Value = $old["K"]   # Supplier Invoice Posting Date
$ws.Range("N8").Value = $old["G"]   # PO Number
$ws.Range("O8").Value = $old["H"]   # Contract Number
$ws.Range("P8").Value = $old["N"]   # Amount
$ws.Range("Q8").Value = $old["O"]   # Supplier Payment Number
$ws.Range("R8").Value = $old["P"]   # Supplier Payment Date
$ws.Range("S8").Value = $old["Q"]   # Due Date
$ws.Range("T8").Value = $old["R"]   # Number of Days
$ws.Range("U8").Value = $old["S"]   # Not Due
$ws.Range("V8").Value = $old["T"]   # Overdue 1-30 Days
$ws.Range("W8").Value = $old["U"]   # Overdue 31-60 Days
$ws.Range("X8").Value = $old["V"]   # Overdue 61-90 Days
$ws.Range("Y8").Value = $old["W"]   # Overdue 91-180 Days

# Two brand-new trailing columns need the same header formatting (yellow,
# bold, centered) as the rest of that band (S8:Y8) - grab it from Y8 (still
# its original style/format) before filling in the new text.
$ws.Range("Y8").Copy() | Out-Null
$ws.Range("Z8:AA8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("Z8").Value = $old["X"]   # Overdue 181-365 Days
$ws.Range("AA8").Value = $old["Y"]  # Overdue over 365 Days

# --- widen the columns whose header text changed to something longer ---
# ColumnWidth is in "characters"; the saved OOXML width is ColumnWidth + 5/6.
$ws.Columns("Q").ColumnWidth = 28.28515625 - (5/6)
$ws.Columns("R").ColumnWidth = 26.7109375 - (5/6)
$ws.Columns("Z").ColumnWidth = 27.85546875 - (5/6)
$ws.Columns("AA").ColumnWidth = 28.140625 - (5/6)

# --- leave the cursor/selection where the author ended up after editing ---
$ws.Range("Q9").Select()
